# Update the "dSF" column (F) values for several rows.
# These correspond to re-pulled / recalculated data per the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -4
    8  = -8
    13 = -8
    14 = -3
    15 = -9
    16 = -4
    17 = -6
    18 = 0
    20 = -2
    21 = 5
    22 = -3
    23 = -3
    25 = -4
    26 = -2
    27 = 11
    28 = 11
    32 = -4
    33 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
